$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "310.93"
    "E2" = "-1.30%"
    "D3" = "48.55"
    "E3" = "8.31%"
    "D4" = "5.240"
    "E4" = "1.88%"
    "D5" = "0.07776"
    "E5" = "-3.72%"
    "D6" = "4.512"
    "E6" = "-0.10%"
    "D7" = "1.293"
    "E7" = "18.34%"
    "D8" = "1.552"
    "E8" = "-8.47%"
    "D9" = "0.1242"
    "E9" = "-4.03%"
    "D10" = "0.1925"
    "E10" = "0.43%"
    "D11" = "0.09301"
    "E11" = "-1.58%"
    "D12" = "0.04555"
    "E12" = "7.04%"
    "E13" = "0.37%"
    "D14" = "0.001291"
    "E14" = "-1.71%"
    "D15" = "0.04206"
    "E15" = "-0.73%"
    "D16" = "0.005851"
    "E16" = "0.34%"
    "D17" = "3.320"
    "E17" = "-2.27%"
    "D18" = "2.405"
    "E18" = "-0.23%"
    "E19" = "2.27%"
    "D20" = "8.131"
    "E20" = "-1.93%"
    "E21" = "-1.09%"
    "D22" = "0.3107"
    "E22" = "-1.25%"
    "D23" = "0.001297"
    "E23" = "1.64%"
    "D24" = "0.004133"
    "E24" = "-3.09%"
    "E25" = "1.18%"
    "D26" = "0.0003568"
    "E26" = "-95.18%"
    "D38" = "0.02573"
    "E38" = "-4.55%"
    "D39" = "0.05789"
    "E39" = "6.16%"
    "D40" = "0.01083"
    "E40" = "97.14%"
    "D41" = "0.007986"
    "E41" = "3.04%"
    "D42" = "0.1427"
    "E42" = "0.24%"
    "D43" = "0.008453"
    "E43" = "14.77%"
    "D44" = "0.008527"
    "E44" = "-0.63%"
    "D45" = "0.3116"
    "E45" = "-0.74%"
    "D46" = "0.00006913"
    "E46" = "1.50%"
    "E47" = "1.19%"
    "D48" = "0.05555"
    "E48" = "-11.20%"
    "D49" = "0.004033"
    "E49" = "1.17%"
    "D50" = "0.00002117"
    "E50" = "1.19%"
    "E51" = "1.19%"
}

foreach ($cell in $updates.Keys) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $updates[$cell]
    $r.ClearFormats()
}
